$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1.xml): update "想去人数" (want-to-go count) values
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 44
$wsExhibit.Range("F4").Value = 249
$wsExhibit.Range("F5").Value = 3901
$wsExhibit.Range("F7").Value = 438

# Sheet "全部类型" (sheet4.xml): update the same events' "想去人数" values
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 44
$wsAll.Range("F4").Value = 249
$wsAll.Range("F5").Value = 3901
$wsAll.Range("F9").Value = 438
